$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing contents (keep formatting) so that shared strings rebuild
# in the exact write order we control below.
$ws.Range("A1:K7").ClearContents()

# --- Header row (A1:K1) - rewritten in original order ---
$ws.Range("A1").Value = "geral_modalidade"
$ws.Range("B1").Value = "mencoes_humor"
$ws.Range("C1").Value = "total"
$ws.Range("D1").Value = "total_sucesso"
$ws.Range("E1").Value = "particip"
$ws.Range("F1").Value = "taxa_sucesso"
$ws.Range("G1").Value = "arrecadado_sucesso"
$ws.Range("H1").Value = "media_sucesso"
$ws.Range("I1").Value = "std_sucesso"
$ws.Range("J1").Value = "min_sucesso"
$ws.Range("K1").Value = "max_sucesso"

# --- New header cells (L1:N1), written right after the existing headers and
# before any data rows, so the new shared strings land right after
# "max_sucesso" and before "aon"/"flex"/"sub" are re-emitted. ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Apply the same header style (bold, border, centered) used by K1 to the new header cells
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (A2:K7) - rewritten in original order; values for E/F are
# rescaled by 100 compared to the source workbook. ---
$ws.Range("A2").Value = "aon"
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 1064
$ws.Range("D2").Value = 633
$ws.Range("E2").Value = 79.70037453183521
$ws.Range("F2").Value = 59.49248120300752
$ws.Range("G2").Value = 17565621.62149074
$ws.Range("H2").Value = 27749.79719034872
$ws.Range("I2").Value = 37329.02017806433
$ws.Range("J2").Value = 41.81688448509265
$ws.Range("K2").Value = 537544.5528256212

$ws.Range("A3").Value = "aon"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = 271
$ws.Range("D3").Value = 197
$ws.Range("E3").Value = 20.29962546816479
$ws.Range("F3").Value = 72.69372693726937
$ws.Range("G3").Value = 6497658.205832288
$ws.Range("H3").Value = 32983.03657782887
$ws.Range("I3").Value = 63531.58925344345
$ws.Range("J3").Value = 54.53892516702949
$ws.Range("K3").Value = 679297.6600721752

$ws.Range("A4").Value = "flex"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 1191
$ws.Range("D4").Value = 1110
$ws.Range("E4").Value = 81.13079019073569
$ws.Range("F4").Value = 93.19899244332494
$ws.Range("G4").Value = 12813838.47206855
$ws.Range("H4").Value = 11543.99862348518
$ws.Range("I4").Value = 23215.13255569839
$ws.Range("J4").Value = 10.77163914429046
$ws.Range("K4").Value = 475290.9541363961

$ws.Range("A5").Value = "flex"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 277
$ws.Range("D5").Value = 273
$ws.Range("E5").Value = 18.86920980926431
$ws.Range("F5").Value = 98.55595667870037
$ws.Range("G5").Value = 5548293.465490552
$ws.Range("H5").Value = 20323.41928751118
$ws.Range("I5").Value = 59929.41679938609
$ws.Range("J5").Value = 76.11778736870863
$ws.Range("K5").Value = 708972.7845446636

$ws.Range("A6").Value = "sub"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 630
$ws.Range("D6").Value = 131
$ws.Range("E6").Value = 92.10526315789474
$ws.Range("F6").Value = 20.79365079365079
$ws.Range("G6").Value = 30146.27685383273
$ws.Range("H6").Value = 230.1242507926162
$ws.Range("I6").Value = 475.50225572397
$ws.Range("J6").Value = 1.087396962410123
$ws.Range("K6").Value = 3475.049171548047

$ws.Range("A7").Value = "sub"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 54
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 7.894736842105263
$ws.Range("F7").Value = 38.88888888888889
$ws.Range("G7").Value = 13040.68090095208
$ws.Range("H7").Value = 620.9848048072417
$ws.Range("I7").Value = 1260.025404745136
$ws.Range("J7").Value = 25.33675915996964
$ws.Range("K7").Value = 5087.076865717208

# --- New data columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes) ---
$ws.Range("L2").Value = 91.11272172566387
$ws.Range("M2").Value = 192790
$ws.Range("N2").Value = 304.565560821485

$ws.Range("L3").Value = 91.82281991764464
$ws.Range("M3").Value = 70763
$ws.Range("N3").Value = 359.2030456852792

$ws.Range("L4").Value = 89.11804758541258
$ws.Range("M4").Value = 143785
$ws.Range("N4").Value = 129.536036036036

$ws.Range("L5").Value = 92.68628097576973
$ws.Range("M5").Value = 59861
$ws.Range("N5").Value = 219.2710622710623

$ws.Range("L6").Value = 18.00852858651895
$ws.Range("M6").Value = 1674
$ws.Range("N6").Value = 12.77862595419847

$ws.Range("L7").Value = 24.42075075084659
$ws.Range("M7").Value = 534
$ws.Range("N7").Value = 25.42857142857143
